# Definition of Done1.1.xlsx - "Add files via upload" edit
# Updates the sign-off table: marks the "Hesse" checkbox column,
# upper-cases the lower-case "x" marks in column J, fixes some text typos
# ("Aufwandbereitschaft" -> "Aufwandsbereitschaft", "Heiser/ Netzler" ->
# "Heiser / Netzler"), adds a new "Dozentenumfrage Aufwandsbereitschaft"
# header in K1, and adds "Hesse " as a new signer in B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (sign-off names) -------------------------------------------
# B8 is a brand new cell -> give it the same look as the other signer
# cells on that row (C8/D8/E8/F8/I8): centered, no border, no fill.
$ws.Range("B8").Value2 = "Hesse "
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").Borders.LineStyle = -4142
$ws.Range("B8").Interior.Pattern = -4142

# --- Row 1 (headers) -------------------------------------------------
# J1: fix typo "Aufwandbereitschaft" -> "Aufwandsbereitschaft"
$ws.Range("J1").Value2 = "E-Mail Dozentenumfrage Aufwandsbereitschaft"

# J8: fix spacing "Heiser/ Netzler" -> "Heiser / Netzler" and align its
# formatting with the rest of row 8 (drops the now-unused "alignment-only"
# style in favour of the shared one already used by C8/D8/E8/F8/I8).
$ws.Range("J8").Value2 = "Heiser / Netzler"
$ws.Range("J8").HorizontalAlignment = -4108
$ws.Range("J8").Borders.LineStyle = -4142
$ws.Range("J8").Interior.Pattern = -4142

# K1: new header cell (same style family as J1 - border + rotated text already
# carried by the column, just set the text)
$ws.Range("K1").Value2 = "Dozentenumfrage Aufwandsbereitschaft"

# --- Rows 2-6: tick the new "Hesse" column (B) and fix the case of the
#     checkmarks in column J (lowercase "x" -> uppercase "X") -----------
$ws.Range("B2").Value2 = "X"
$ws.Range("B3").Value2 = "X"
$ws.Range("B4").Value2 = "X"
$ws.Range("B5").Value2 = "X"
$ws.Range("B6").Value2 = "X"

$ws.Range("J2").Value2 = "X"
$ws.Range("J3").Value2 = "X"
$ws.Range("J4").Value2 = "X"
$ws.Range("J5").Value2 = "X"
$ws.Range("J6").Value2 = "X"

# --- Column widths: accommodate the new / changed content --------------
# (target widths match Excel's native autofit grid; values below are the
# closest achievable ColumnWidth inputs on that grid)
$ws.Columns.Item(1).ColumnWidth = 33.833333333333336
$ws.Columns.Item(2).ColumnWidth = 5.666666666666667
$ws.Columns.Item(3).ColumnWidth = 6.333333333333334
$ws.Columns.Item(4).ColumnWidth = 6.833333333333334
$ws.Columns.Item(5).ColumnWidth = 6
$ws.Columns.Item(6).ColumnWidth = 5.833333333333334
$ws.Columns.Item(7).ColumnWidth = 3.8333333333333335
$ws.Columns.Item(8).ColumnWidth = 3.8333333333333335
$ws.Columns.Item(9).ColumnWidth = 12.333333333333334
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666
$ws.Columns.Item(13).ColumnWidth = 11.5

# --- Row height for the header row (minor re-wrap adjustment) ----------
$ws.Rows.Item(1).RowHeight = 178.9

# --- Selection cursor, matching the saved selection in the workbook -----
$ws.Range("K2").Select()
